$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 6: status flips from "done" to "open"; the date + tech cells are cleared ---
$ws.Range("E6").Value = "open"
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()

# --- Insert the first block of 3 new rows at row 8 (pagination + two blank placeholder rows) ---
$ws.Rows("8:10").Insert()

$ws.Range("C8").Value = "pagination"
$ws.Range("E8").Value = "done"
$ws.Range("F8").Value = 43988
$ws.Range("G8").Value = "ajax "

# --- Insert the second block of 2 new rows at row 14 (browse / grid view feature) ---
$ws.Rows("14:15").Insert()

$ws.Range("C15").Value = "cart items in sidebarB"
$ws.Range("E15").Value = "done"
$ws.Range("F15").NumberFormat = $ws.Range("F3").NumberFormat
$ws.Range("F15").Value = 43989

$ws.Range("B14").Value = "browse"
$ws.Range("C14").Value = "grid view/list view"

# --- Insert the third block of 2 new rows at row 19 (cart totalprice / cart badge) ---
$ws.Rows("19:20").Insert()

$ws.Range("C19").Value = "cart totalprice in sidebarB"
$ws.Range("E19").Value = "done"
$ws.Range("F19").NumberFormat = $ws.Range("F3").NumberFormat
$ws.Range("F19").Value = 43989

$ws.Range("C20").Value = "cart badge with quantity"
$ws.Range("E20").Value = "done"
$ws.Range("F20").NumberFormat = $ws.Range("F3").NumberFormat
$ws.Range("F20").Value = 43989

# --- Match the final selection + print orientation left behind by the edit ---
$ws.PageSetup.Orientation = 1
$ws.Range("E20:F20").Select()
